# Pop the next queued name off "Sheet1" (row 2 -> "fuxp1un3") and record it
# as consumed on the "used" log sheet, together with the source filename and
# the timestamp it was used at.

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# The value being popped (row 2 of Sheet1, right under the header in row 1).
$poppedId = $namesSheet.Range("A2").Value()

# Remove that row entirely; rows below shift up by one (474 -> 473 rows),
# shrinking Sheet1's used range from A1:A474 to A1:A473.
$namesSheet.Range("A2").EntireRow.Delete()

# Append the popped id as a new row on the "used" sheet, growing its used
# range from A1:C25 to A1:C26.
$nextRow = $usedSheet.UsedRange.Rows.Count + 1
$usedSheet.Cells.Item($nextRow, 1).Value = $poppedId
$usedSheet.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月21日 15_24_19.png"
$usedSheet.Cells.Item($nextRow, 3).Value = "2026-01-21 15:26:11"
